$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text before writing, so numeric-looking price strings
# (e.g. "58.16") are preserved verbatim as text instead of being parsed as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '70.091.23'
$ws.Range("E2").Value = '  +7.19%  '
# Row 3
$ws.Range("D3").Value = '3.597.58'
$ws.Range("E3").Value = '  +6.67%  '
# Row 4
$ws.Range("E4").Value = '  -0.11%  '
# Row 5
$ws.Range("D5").Value = '594.55'
$ws.Range("E5").Value = '  +6.39%  '
# Row 6
$ws.Range("D6").Value = '192.82'
$ws.Range("E6").Value = '  +10.63%  '
# Row 7
$ws.Range("D7").Value = '0.649'
$ws.Range("E7").Value = '  +3.25%  '
# Row 8
$ws.Range("D8").Value = '3.592.92'
$ws.Range("E8").Value = '  +6.92%  '
# Row 9
$ws.Range("D9").Value = '1.00'
$ws.Range("E9").Value = '  -0.12%  '
# Row 10
$ws.Range("E10").Value = '  +5.13%  '
# Row 11
$ws.Range("D11").Value = '0.664'
$ws.Range("E11").Value = '  +5.18%  '
# Row 12
$ws.Range("D12").Value = '58.16'
$ws.Range("E12").Value = '  +9.42%  '
# Row 13
$ws.Range("D13").Value = '0.0000296'
$ws.Range("E13").Value = '  +7.27%  '
# Row 14
$ws.Range("D14").Value = '9.76'
$ws.Range("E14").Value = '  +6.60%  '
# Row 15
$ws.Range("D15").Value = '4.166.35'
$ws.Range("E15").Value = '  +5.92%  '
# Row 16
$ws.Range("E16").Value = '  +6.31%  '
# Row 17
$ws.Range("D17").Value = '3.595.59'
$ws.Range("E17").Value = '  +5.76%  '
# Row 18
$ws.Range("D18").Value = '69.972.35'
$ws.Range("E18").Value = '  +6.88%  '
# Row 19
$ws.Range("D19").Value = '12.66'
$ws.Range("E19").Value = '  +7.42%  '
# Row 20
$ws.Range("E20").Value = '  +1.27%  '
# Row 21
$ws.Range("E21").Value = '  +5.93%  '
# Row 22
$ws.Range("D22").Value = '499.41'
$ws.Range("E22").Value = '  +4.18%  '
# Row 23
$ws.Range("D23").Value = '5.50'
$ws.Range("E23").Value = '  +10.41%  '
# Row 24
$ws.Range("D24").Value = '17.00'
$ws.Range("E24").Value = '  +19.00%  '
# Row 25
$ws.Range("D25").Value = '4.48'
$ws.Range("E25").Value = '  +9.54%  '
# Row 26
$ws.Range("D26").Value = '91.20'
$ws.Range("E26").Value = '  +1.24%  '
# Row 27
$ws.Range("E27").Value = '  +6.97%  '
# Row 28
$ws.Range("D28").Value = '11.25'
$ws.Range("E28").Value = '  +6.66%  '
# Row 29
$ws.Range("D29").Value = '9.37'
$ws.Range("E29").Value = '  +7.96%  '
# Row 30
$ws.Range("D30").Value = '32.36'
$ws.Range("E30").Value = '  +4.05%  '
# Row 31
$ws.Range("D31").Value = '7.54'
$ws.Range("E31").Value = '  +15.59%  '
# Row 32
$ws.Range("D32").Value = '12.21'
$ws.Range("E32").Value = '  +7.22%  '
# Row 33
$ws.Range("D33").Value = '616.65'
$ws.Range("E33").Value = '  +7.76%  '
# Row 34
$ws.Range("D34").Value = '65.39'
$ws.Range("E34").Value = '  +2.70%  '
# Row 35
$ws.Range("E35").Value = '  +8.25%  '
# Row 36
$ws.Range("E36").Value = '  +13.65%  '
# Row 37
$ws.Range("D37").Value = '0.148'
$ws.Range("E37").Value = '  +4.62%  '
# Row 38
$ws.Range("B38").Value = 'Stacks'
$ws.Range("C38").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D38").Value = '3.72'
$ws.Range("E38").Value = '  +2.88%  '
# Row 39
$ws.Range("E39").Value = '  -0.02%  '
# Row 40
$ws.Range("B40").Value = 'InjectiveProtocol'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D40").Value = '38.04'
$ws.Range("E40").Value = '  +7.01%  '
# Row 41
$ws.Range("E41").Value = '  +7.42%  '
# Row 42
$ws.Range("D42").Value = '3.340.49'
$ws.Range("E42").Value = '  +7.74%  '
# Row 43
$ws.Range("D43").Value = '3.11'
$ws.Range("E43").Value = '  +11.43%  '
# Row 44
$ws.Range("E44").Value = '  +11.23%  '
# Row 45
$ws.Range("E45").Value = '  +7.32%  '
# Row 46
$ws.Range("E46").Value = '  +16.71%  '
# Row 47
$ws.Range("B47").Value = 'Stellar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D47").Value = '0.138'
$ws.Range("E47").Value = '  +3.26%  '
# Row 48
$ws.Range("B48").Value = 'ApeXProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D48").Value = '3.26'
$ws.Range("E48").Value = '  +3.32%  '
# Row 49
$ws.Range("D49").Value = '9.05'
$ws.Range("E49").Value = '  +7.82%  '
# Row 50
$ws.Range("B50").Value = 'LidoDAOToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D50").Value = '3.26'
$ws.Range("E50").Value = '  +6.18%  '
# Row 51
$ws.Range("B51").Value = 'FirstDigitalUSD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D51").Value = '0.999'
$ws.Range("E51").Value = '  -0.13%  '

# Restore column D to the default (unstyled) cell style; only the
# number-format override was needed to keep the values textual.
$ws.Range("D2:D51").Style = "Normal"
